$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$titleShape = $s.Shapes.Item(3)

# Update position/size of the title placeholder (values are in points; 1 pt = 12700 EMU)
$titleShape.Left = 0
$titleShape.Top = -0.00007874015748031496
$titleShape.Width = 720
$titleShape.Height = 73.687485

# Update the text content
$titleShape.TextFrame.TextRange.Text = "Which GICS sector spends the most in R&D, taking into account the Operating Income?"

# Center-align the paragraph
$titleShape.TextFrame.TextRange.ParagraphFormat.Alignment = 2
